# Updated cryptos list on Thu Nov  9 10:26:50 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "36.931.61"
$ws.Range("E2").Value = "  +4.31%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.917.98"
$ws.Range("E3").Value = "  +1.61%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - was BNB, now XRP
$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.707"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.96%  "

# Row 6 - was XRP, now BNB
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "251.22"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.82%  "

# Row 7 - USDC
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8 - Solana
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "47.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +8.57%  "

# Row 9 - Cardano
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.374"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.59%  "

# Row 10 - OKB
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "58.53"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +9.48%  "

# Row 11 - Dogecoin
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0765"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.80%  "

# Row 12 - TRON
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0999"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.09%  "

# Row 13 - Chainlink
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "14.68"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +8.21%  "

# Row 14 - Polygon
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.818"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +5.67%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.195.39"
$ws.Range("E15").Value = "  +1.65%  "

# Row 16 - Polkadot
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.16"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +4.02%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.914.77"
$ws.Range("E17").Value = "  +1.25%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "37.278.67"
$ws.Range("E18").Value = "  +5.27%  "

# Row 19 - Litecoin
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "74.99"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.01%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0864"
$ws.Range("E20").Value = "  +4.27%  "

# Row 21 - was BitcoinCash, now Avalanche
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.68"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +6.31%  "

# Row 22 - was Avalanche, now BitcoinCash
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "252.40"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.86%  "

# Row 23 - Uniswap
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.21"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24 - Toncoin
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.65"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.90%  "

# Row 25 - Dai
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "

# Row 26 - PancakeSwap
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.22"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.41%  "

# Row 27 - Monero
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "168.06"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.86%  "

# Row 28 - Cosmos
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.85"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.04%  "

# Row 29 - EthereumClassic
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "18.80"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.39%  "

# Row 30 - Stellar
$ws.Range("E30").Value = "  +1.84%  "

# Row 31 - Filecoin
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.60"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +6.69%  "

# Row 32 - Hedera
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0621"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +4.24%  "

# Row 33 - Kaspa
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0912"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +24.10%  "

# Row 35 - WEMIXToken
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.89"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.53%  "

# Row 36 - BinanceUSD
$ws.Range("E36").Value = "  +0.06%  "

# Row 37 - TrustWalletToken
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.52"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.58%  "

# Row 38 - Gas
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "18.57"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +55.96%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +2.15%  "

# Row 40 - LidoDAOToken
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.02"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.01%  "

# Row 41 - Aave
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "105.86"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +8.49%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  +4.53%  "

# Row 43 - InjectiveProtocol
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "17.96"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.64%  "

# Row 44 - HuobiToken
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.83"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +18.14%  "

# Row 45 - ARBITRUM
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.11"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.62%  "

# Row 46 - Maker
$ws.Range("D46").Value = "1.353.20"
$ws.Range("E46").Value = "  +3.18%  "

# Row 47 - RenderToken
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.40"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.96%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +1.07%  "

# Row 49 - MXToken
$ws.Range("E49").Value = "  +2.17%  "

# Row 50 - FraxShare
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.49"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.26%  "

# Row 51 - MultiversX
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "43.51"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.92%  "

